# "Generate Report for Archive"
#
# The localization-status report is regenerated: the row that used to
# describe 603796dc-fce3-494f-9d55-e04697052c8e.md and the row describing
# 8ba283ea-b616-46c2-913d-2969658428ec.md trade places (the 8ba283ea entry
# now shows up before the 603796dc entry, and its status moved on from
# "Ready for handoff" to "In Translation"). This touches the Overview sheet
# as well as the per-locale (zh-cn / de-de) detail sheets, including the
# hyperlinks that decorate the file-name cells.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A5").Value = "8ba283ea-b616-46c2-913d-2969658428ec.md"
$ws.Range("B5").Value = "e2e\8ba283ea-b616-46c2-913d-2969658428ec.md"
$ws.Range("E5").Value = "In Translation"
$ws.Range("F5").Value = "In Translation"
$ws.Range("G5").Value = "2016-08-17 04:37:49"

$ws.Range("A6").Value = "603796dc-fce3-494f-9d55-e04697052c8e.md"
$ws.Range("B6").Value = "e2e\603796dc-fce3-494f-9d55-e04697052c8e.md"
$ws.Range("E6").Value = "Ready for handoff"
$ws.Range("F6").Value = "Ready for handoff"
$ws.Range("G6").Value = "2016-08-17 04:35:39"

# Rebuild the hyperlinks on column B so the display text follows the new
# row content while each row keeps pointing at the commit it always did.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/097dc7390381a7e240c243e97ca7109f7b3ac8a7/e2e/b9cccdbb-31bc-4522-bc8d-ac4d8d4f85ea.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "e2e\b9cccdbb-31bc-4522-bc8d-ac4d8d4f85ea.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/173e8d20cf4a7687fd3edf86d9afdbe392de6000/e2e/00ba7319-81d3-43ad-aa3f-e688dcfce039.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "e2e\00ba7319-81d3-43ad-aa3f-e688dcfce039.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/173e8d20cf4a7687fd3edf86d9afdbe392de6000/e2e/5d9e857b-b1fc-4d60-aade-345b6b086975.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "e2e\5d9e857b-b1fc-4d60-aade-345b6b086975.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f16fdbaa46a1b729355cb025932a4faa4d43c8c2/e2e/603796dc-fce3-494f-9d55-e04697052c8e.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "e2e\8ba283ea-b616-46c2-913d-2969658428ec.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8b12eb8a5cd06b0dab8ed7dcdfecaf3b98f9ce6a/e2e/8ba283ea-b616-46c2-913d-2969658428ec.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "e2e\603796dc-fce3-494f-9d55-e04697052c8e.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5e2912e2cb46f3d9b6e92ccf878cf65f2475c127/e2e/bafcc2e9-6684-44b8-87dc-675f32e05915.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "e2e\bafcc2e9-6684-44b8-87dc-675f32e05915.md") | Out-Null

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A5").Value = "8ba283ea-b616-46c2-913d-2969658428ec.md"
$ws.Range("C5").Value = "In Translation"
$ws.Range("G5").Value = "8ba283ea-b616-46c2-913d-2969658428ec.9398e1d4a6c32d882cbfe639bf55260f2341bae7.zh-cn.xlf"
$ws.Range("H5").Value = "2016-08-17 04:37:44"

$ws.Range("A6").Value = "603796dc-fce3-494f-9d55-e04697052c8e.md"
$ws.Range("C6").Value = "Ready for handoff"
$ws.Range("G6").Value = "603796dc-fce3-494f-9d55-e04697052c8e.dfb865b02a6851becd229810cfb583045bc5f519.zh-cn.xlf"
$ws.Range("H6").Value = "2016-08-17 04:35:35"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/097dc7390381a7e240c243e97ca7109f7b3ac8a7/e2e/b9cccdbb-31bc-4522-bc8d-ac4d8d4f85ea.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "b9cccdbb-31bc-4522-bc8d-ac4d8d4f85ea.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/a34d3fe2585dac1ef0051a911e5de0f75ae0a3a3/e2e/b9cccdbb-31bc-4522-bc8d-ac4d8d4f85ea.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "b9cccdbb-31bc-4522-bc8d-ac4d8d4f85ea.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/173e8d20cf4a7687fd3edf86d9afdbe392de6000/e2e/00ba7319-81d3-43ad-aa3f-e688dcfce039.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "00ba7319-81d3-43ad-aa3f-e688dcfce039.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/173e8d20cf4a7687fd3edf86d9afdbe392de6000/e2e/5d9e857b-b1fc-4d60-aade-345b6b086975.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "5d9e857b-b1fc-4d60-aade-345b6b086975.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f16fdbaa46a1b729355cb025932a4faa4d43c8c2/e2e/603796dc-fce3-494f-9d55-e04697052c8e.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "8ba283ea-b616-46c2-913d-2969658428ec.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8b12eb8a5cd06b0dab8ed7dcdfecaf3b98f9ce6a/e2e/8ba283ea-b616-46c2-913d-2969658428ec.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "603796dc-fce3-494f-9d55-e04697052c8e.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5e2912e2cb46f3d9b6e92ccf878cf65f2475c127/e2e/bafcc2e9-6684-44b8-87dc-675f32e05915.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "bafcc2e9-6684-44b8-87dc-675f32e05915.md") | Out-Null

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A5").Value = "8ba283ea-b616-46c2-913d-2969658428ec.md"
$ws.Range("C5").Value = "In Translation"
$ws.Range("G5").Value = "8ba283ea-b616-46c2-913d-2969658428ec.9398e1d4a6c32d882cbfe639bf55260f2341bae7.de-de.xlf"
$ws.Range("H5").Value = "2016-08-17 04:37:49"

$ws.Range("A6").Value = "603796dc-fce3-494f-9d55-e04697052c8e.md"
$ws.Range("C6").Value = "Ready for handoff"
$ws.Range("G6").Value = "603796dc-fce3-494f-9d55-e04697052c8e.dfb865b02a6851becd229810cfb583045bc5f519.de-de.xlf"
$ws.Range("H6").Value = "2016-08-17 04:35:39"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/097dc7390381a7e240c243e97ca7109f7b3ac8a7/e2e/b9cccdbb-31bc-4522-bc8d-ac4d8d4f85ea.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "b9cccdbb-31bc-4522-bc8d-ac4d8d4f85ea.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/86631f3cf1b454092da7c50bd512bc5c471ba916/e2e/b9cccdbb-31bc-4522-bc8d-ac4d8d4f85ea.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "b9cccdbb-31bc-4522-bc8d-ac4d8d4f85ea.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/173e8d20cf4a7687fd3edf86d9afdbe392de6000/e2e/00ba7319-81d3-43ad-aa3f-e688dcfce039.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "00ba7319-81d3-43ad-aa3f-e688dcfce039.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/173e8d20cf4a7687fd3edf86d9afdbe392de6000/e2e/5d9e857b-b1fc-4d60-aade-345b6b086975.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "5d9e857b-b1fc-4d60-aade-345b6b086975.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f16fdbaa46a1b729355cb025932a4faa4d43c8c2/e2e/603796dc-fce3-494f-9d55-e04697052c8e.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "8ba283ea-b616-46c2-913d-2969658428ec.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8b12eb8a5cd06b0dab8ed7dcdfecaf3b98f9ce6a/e2e/8ba283ea-b616-46c2-913d-2969658428ec.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "603796dc-fce3-494f-9d55-e04697052c8e.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5e2912e2cb46f3d9b6e92ccf878cf65f2475c127/e2e/bafcc2e9-6684-44b8-87dc-675f32e05915.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "bafcc2e9-6684-44b8-87dc-675f32e05915.md") | Out-Null
